{"js": "// The document contains two \"<id>...</id>\" tags (one for the page/article\n// id, one for the figure id further down). We must only touch the first\n// one, which currently reads \"<id>p080v_a1</id>\" split across three runs:\n//   1) \"<id>\"       (Courier New, color 7f6000, size 9pt)\n//   2) \"p080v_a1\"   (color 000000)\n//   3) \"</id>\"      (Courier New, color 7f6000, size 9pt)\n// The edit merges these three runs into a single run reading\n// \"<id>p080v_1</id>\" that keeps the formatting of the first run.\n\nconst body = context.document.body;\n\n// \"p080v_a1\" is unique in the document, so it safely anchors the middle run.\nconst middleMatches = body.search(\"p080v_a1\", { matchCase: true, matchWholeWord: false });\nmiddleMatches.load(\"text\");\nawait context.sync();\n\nif (middleMatches.items.length === 0) {\n  throw new Error(\"Could not find the 'p080v_a1' run to edit.\");\n}\nconst middleRange = middleMatches.items[0];\n\n// There are multiple \"<id>\" / \"</id>\" occurrences in the document (the\n// figure's id tag lower down uses the same literal text), so disambiguate\n// by picking the \"<id>\" immediately before, and the \"</id>\" immediately\n// after, our middle run.\nconst openMatches = body.search(\"<id>\", { matchCase: true });\nopenMatches.load(\"text\");\nconst closeMatches = body.search(\"</id>\", { matchCase: true });\ncloseMatches.load(\"text\");\nawait context.sync();\n\nlet openRange = null;\nfor (let i = 0; i < openMatches.items.length; i++) {\n  const cmp = openMatches.items[i].compareLocationWith(middleRange);\n  await context.sync();\n  if (cmp.value === \"AdjacentBefore\") {\n    openRange = openMatches.items[i];\n    break;\n  }\n}\n\nlet closeRange = null;\nfor (let i = 0; i < closeMatches.items.length; i++) {\n  const cmp = middleRange.compareLocationWith(closeMatches.items[i]);\n  await context.sync();\n  if (cmp.value === \"AdjacentBefore\") {\n    closeRange = closeMatches.items[i];\n    break;\n  }\n}\n\nif (!openRange || !closeRange) {\n  throw new Error(\"Could not locate the surrounding <id>/</id> runs.\");\n}\n\n// Combine the three runs (\"<id>\", \"p080v_a1\", \"</id>\") into a single range\n// and replace its text. The new run inherits the formatting of the first\n// run in the combined range (the \"<id>\" run's Courier New / 7f6000 / 9pt).\nconst fullRange = openRange.expandTo(closeRange);\nfullRange.insertText(\"<id>p080v_1</id>\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The document contains two \"<id>...</id>\" tags (one for the page/article\n# id near the top, one for a figure further down). Only the first one,\n# currently \"<id>p080v_a1</id>\" split across three runs:\n#   1) \"<id>\"       (Courier New, color 7f6000, size 9pt)\n#   2) \"p080v_a1\"   (color 000000)\n#   3) \"</id>\"      (Courier New, color 7f6000, size 9pt)\n# must be merged into a single run reading \"<id>p080v_1</id>\" that keeps the\n# formatting of the first (\"<id>\") run.\n\n$d = $word.ActiveDocument\n\n# \"p080v_a1\" is unique in the document, so Find.Execute unambiguously\n# anchors the middle run (unlike bare \"<id>\"/\"</id>\", which also occur\n# around the unrelated figure id lower in the document). Find.Execute\n# collapses $middle to the matched text's own boundaries.\n$middle = $d.Content\n$find = $middle.Find\n$find.Text = \"p080v_a1\"\n$find.MatchCase = $true\n$find.Forward = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the 'p080v_a1' run to edit.\"\n}\n\n# Grow the range by the 4 characters of the literal \"<id>\" before it and\n# the 5 characters of the literal \"</id>\" after it, combining all three\n# runs into one range.\n$fullRange = $d.Range($middle.Start - 4, $middle.End + 5)\n\n# Sanity-check we grabbed exactly the tag we expect before mutating it.\nif ($fullRange.Text -ne \"<id>p080v_a1</id>\") {\n    throw \"Unexpected surrounding text: [$($fullRange.Text)]\"\n}\n\n# Replacing .Text on the combined range merges it into a single run that\n# inherits the formatting already present at the start of the range (the\n# \"<id>\" run's Courier New / 7f6000 / 9pt), matching the target edit.\n$fullRange.Text = \"<id>p080v_1</id>\"\n"}
